$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) are stored as plain text in the sheet
# (numeric- and percent-looking strings). Force the Text number format on the
# whole touched range first so Excel does not reinterpret these as numbers/
# percentages when we assign the new string values, then restore the default
# style so no stray per-cell formatting is left behind.
$priceVolRange = $ws.Range("D2:E47")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "286.30"
$ws.Range("E2").Value = "1.73%"
$ws.Range("E3").Value = "3.09%"
$ws.Range("D4").Value = "5.074"
$ws.Range("E4").Value = "0.63%"
$ws.Range("D5").Value = "0.06749"
$ws.Range("E5").Value = "4.35%"
$ws.Range("D6").Value = "7.348"
$ws.Range("E6").Value = "1.90%"
$ws.Range("D7").Value = "1.379"
$ws.Range("E7").Value = "-0.79%"
$ws.Range("D8").Value = "0.9007"
$ws.Range("E8").Value = "-3.06%"
$ws.Range("D9").Value = "0.1594"
$ws.Range("E9").Value = "3.72%"
$ws.Range("D10").Value = "0.06870"
$ws.Range("E10").Value = "9.41%"
$ws.Range("D11").Value = "0.07567"
$ws.Range("E11").Value = "-0.07%"
$ws.Range("D12").Value = "0.02920"
$ws.Range("E12").Value = "1.04%"
$ws.Range("D13").Value = "0.08999"
$ws.Range("E13").Value = "0.19%"
$ws.Range("D14").Value = "0.001565"
$ws.Range("E14").Value = "-1.54%"
$ws.Range("D15").Value = "0.0006475"
$ws.Range("E15").Value = "0.81%"
$ws.Range("D16").Value = "0.006565"
$ws.Range("E16").Value = "7.61%"
$ws.Range("D17").Value = "3.457"
$ws.Range("E17").Value = "0.48%"
$ws.Range("D18").Value = "3.437"
$ws.Range("E18").Value = "1.62%"
$ws.Range("D19").Value = "2.231"
$ws.Range("E19").Value = "-0.14%"
$ws.Range("D20").Value = "0.3206"
$ws.Range("E20").Value = "0.69%"
$ws.Range("D21").Value = "0.1320"
$ws.Range("E21").Value = "2.99%"
$ws.Range("D22").Value = "4.001"
$ws.Range("E22").Value = "-1.30%"
$ws.Range("D23").Value = "0.1556"
$ws.Range("E23").Value = "0.26%"
$ws.Range("D24").Value = "0.04484"
$ws.Range("E24").Value = "2.01%"
$ws.Range("D25").Value = "0.001202"
$ws.Range("E25").Value = "0.84%"
$ws.Range("D26").Value = "0.004379"
$ws.Range("E26").Value = "-0.10%"
$ws.Range("E27").Value = "-7.21%"
$ws.Range("D28").Value = "0.0001614"
$ws.Range("E28").Value = "-0.84%"
$ws.Range("D40").Value = "0.04247"
$ws.Range("E40").Value = "3.73%"
$ws.Range("D41").Value = "0.006809"
$ws.Range("E41").Value = "1.09%"
$ws.Range("D42").Value = "0.1241"
$ws.Range("E42").Value = "1.70%"
$ws.Range("D43").Value = "0.002194"
$ws.Range("E43").Value = "3.85%"
$ws.Range("D44").Value = "0.01140"
$ws.Range("E44").Value = "-5.74%"
$ws.Range("D45").Value = "0.00005719"
$ws.Range("E45").Value = "0.88%"
$ws.Range("D46").Value = "1.929"
$ws.Range("E46").Value = "-1.85%"
$ws.Range("D47").Value = "0.01304"
$ws.Range("E47").Value = "-0.31%"

$priceVolRange.Style = "Normal"

# Columns B (Coin) and C (Link) are plain text already; no special handling needed.
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
